# Apply cryptos list update (Mon Jun 10 13:59:22 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.365.51"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "3.671.72"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'644.84"
$ws.Range("E5").Value = "  -5.21%  "
$ws.Range("D6").Value = "'158.56"
$ws.Range("E6").Value = "  -1.13%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("E10").Value = "  -1.38%  "
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("D13").Value = "4.290.43"
$ws.Range("E13").Value = "  -0.44%  "
$ws.Range("D14").Value = "'32.46"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("D15").Value = "3.717.09"
$ws.Range("E15").Value = "  +1.03%  "
$ws.Range("D16").Value = "69.338.21"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").Value = "'15.88"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "'6.45"
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").Value = "'465.53"
$ws.Range("E20").Value = "  -1.39%  "
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").Value = "'0.645"
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("D23").Value = "'79.28"
$ws.Range("E23").Value = "  -1.08%  "
$ws.Range("D24").Value = "3.818.48"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -0.22%  "
$ws.Range("D27").Value = "'10.69"
$ws.Range("D28").Value = "'8.97"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("D31").Value = "'2.00"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("D33").Value = "'26.87"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").Value = "3.663.74"
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  -0.35%  "
$ws.Range("D37").Value = "'8.41"
$ws.Range("E37").Value = "  +0.23%  "
$ws.Range("D39").Value = "'178.71"
$ws.Range("E39").Value = "  +5.50%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "'5.83"
$ws.Range("E40").Value = "  -6.46%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'2.20"
$ws.Range("E42").Value = "  -2.68%  "
$ws.Range("D43").Value = "'0.0891"
$ws.Range("E43").Value = "  -1.70%  "
$ws.Range("D44").Value = "'0.924"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "'28.20"
$ws.Range("E46").Value = "  -1.30%  "
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("E48").Value = "  -2.41%  "
$ws.Range("E49").Value = "  -5.01%  "
$ws.Range("D50").Value = "'7.75"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("E51").Value = "  -4.06%  "
